$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-04-04 13:23:06"

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = $newTimestamp
}
